$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.130.35'
$ws.Range("E2").Value = '  +1.27%  '

$ws.Range("D3").Value = '2.649.66'
$ws.Range("E3").Value = '  -0.24%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '607.81'
$ws.Range("E5").Value = '  -0.36%  '

$ws.Range("D6").Value = '148.64'
$ws.Range("E6").Value = '  +3.28%  '

$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '0.590'
$ws.Range("E8").Value = '  +0.67%  '

$ws.Range("E9").Value = '  +1.59%  '

$ws.Range("D10").Value = '0.388'
$ws.Range("E10").Value = '  +7.58%  '

$ws.Range("D11").Value = '5.62'
$ws.Range("E11").Value = '  -0.05%  '

$ws.Range("E12").Value = '  -0.77%  '

$ws.Range("D13").Value = '27.61'
$ws.Range("E13").Value = '  +1.07%  '

$ws.Range("D14").Value = '3.120.30'
$ws.Range("E14").Value = '  -0.27%  '

$ws.Range("D15").Value = '63.967.41'
$ws.Range("E15").Value = '  +1.26%  '

$ws.Range("E16").Value = '  +1.82%  '

$ws.Range("D17").Value = '2.653.43'
$ws.Range("E17").Value = '  -1.36%  '

$ws.Range("D18").Value = '11.98'
$ws.Range("E18").Value = '  +4.81%  '

$ws.Range("D19").Value = '4.61'
$ws.Range("E19").Value = '  +4.30%  '

$ws.Range("D20").Value = '347.30'
$ws.Range("E20").Value = '  +1.72%  '

$ws.Range("D21").Value = '6.91'
$ws.Range("E21").Value = '  +0.68%  '

$ws.Range("E22").Value = '  +0.14%  '

$ws.Range("D23").Value = '5.56'
$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("D24").Value = '66.26'
$ws.Range("E24").Value = '  -0.98%  '

$ws.Range("D25").Value = '1.70'
$ws.Range("E25").Value = '  +9.58%  '

$ws.Range("E26").Value = '  +4.12%  '

$ws.Range("D27").Value = '9.34'
$ws.Range("E27").Value = '  +7.63%  '

$ws.Range("D28").Value = '557.02'
$ws.Range("E28").Value = '  +2.02%  '

$ws.Range("D29").Value = '8.17'
$ws.Range("E29").Value = '  +4.50%  '

$ws.Range("E30").Value = '  -1.30%  '

$ws.Range("E31").Value = '  -0.02%  '

$ws.Range("E32").Value = '  +0.93%  '

$ws.Range("D33").Value = '0.0₃0850'
$ws.Range("E33").Value = '  +5.30%  '

$ws.Range("E34").Value = '  -0.58%  '

$ws.Range("E35").Value = '  +4.18%  '

$ws.Range("D36").Value = '168.33'
$ws.Range("E36").Value = '  -2.18%  '

$ws.Range("D37").Value = '0.407'
$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.08%  '

$ws.Range("E39").Value = '  +3.95%  '

$ws.Range("D40").Value = '19.30'
$ws.Range("E40").Value = '  +0.98%  '

$ws.Range("E41").Value = '  +0.02%  '

$ws.Range("D42").Value = '167.68'
$ws.Range("E42").Value = '  -3.80%  '

$ws.Range("D43").Value = '40.30'
$ws.Range("E43").Value = '  +0.54%  '

$ws.Range("D44").Value = '3.84'
$ws.Range("E44").Value = '  +2.44%  '

$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").Value = '0.0571'
$ws.Range("E45").Value = '  -0.42%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '21.98'
$ws.Range("E46").Value = '  -1.06%  '

$ws.Range("E47").Value = '  -0.37%  '

$ws.Range("E48").Value = '  +2.25%  '

$ws.Range("E49").Value = '  +13.34%  '

$ws.Range("E50").Value = '  +0.05%  '

$ws.Range("D51").Value = '19.02'
$ws.Range("E51").Value = '  +1.58%  '

